$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Battery holder update ---
# Fix footprint name (remove "*** NEED TO EDIT***")
$ws.Range("E8").Value = "Battery:BatteryHolder_Keystone_1060_1x2032"

# Swap the product URL to the cheaper 1060TR part and update price
$ws.Range("H8").Value = "https://www.digikey.ca/en/products/detail/keystone-electronics/1060TR/303557"
$ws.Range("I8").Value = 2.12

# Add a real hyperlink on the datasheet cell (text stays the same)
$ws.Hyperlinks.Add($ws.Range("G8"), $ws.Range("G8").Value) | Out-Null

# --- Row 16: add related component link ---
$ws.Range("J16").Value = "https://www.digikey.ca/en/products/detail/vishay-general-semiconductor-diodes-division/BZX884B5V6L-G3-08/14312759"

# --- Row 24: add real hyperlink on inductor datasheet cell ---
$ws.Hyperlinks.Add($ws.Range("G24"), $ws.Range("G24").Value) | Out-Null

# --- Row 25: clear the stray "m" note in Column1 ---
$ws.Range("M25").Value = $null

# --- New row 42 ---
$ws.Range("D42").Value = "e"
$ws.Range("E42").Value = "e"
$ws.Range("L42").Formula = "=Table1[[#This Row],[Qnty]]*(Table1[[#This Row],[Price]]+Table1[[#This Row],[Price2]])"

# --- Expand Table1 to include the new row ---
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A7:M42"))

# --- Update total formula to cover the new row ---
$ws.Range("B6").Formula = "=SUM(L8:L42)"

$wb.Save()
